$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new row of data (row 8), mirroring the layout of rows 6/7
# A8: same text as A6/A7 ("2/13/2020jaclemon")
# C8: new duration string "15 minutes"
# D8: new description string
$ws.Range("A8").Value = "2/13/2020jaclemon"
$ws.Range("C8").Value = "15 minutes"
$ws.Range("D8").Value = "Used Clion to prevent from repeated command lines from having an affect"

$ws.Range("E13").Select()
